# skill_prereq sheet gets two new skill types in the dropdown lookup table
# (FAITH, PSIONICS, SPELLCASTING, WEIRD_SCIENCE) plus a new "school"
# column (D) used by a handful of edges (Battle / Arcana / ANY), and a
# large batch of new skill-prerequisite rows for edges that previously had
# no skill requirement recorded. The old sparse edge list (rows 2-113,
# one row per edge, most blank) is replaced by a compact list containing
# only the edges that actually carry a skill prerequisite, and the lookup
# tables (rows 119-141 / 119-123) move up to rows 50-76 / 50-54.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("skill_prereq")

# --- data -------------------------------------------------------------

$edgeRows = @(
    @("BLOCK","FIGHTING","D8",""),
    @("COUNTERATTACK","FIGHTING","D8",""),
    @("FLORENTINE","FIGHTING","D8",""),
    @("FRENZY","FIGHTING","D10",""),
    @("MARTIAL_ARTIST","FIGHTING","D6",""),
    @("IMPROVED_MARTIAL_ARTIST","FIGHTING","D10",""),
    @("ROCK_AND_ROLL","SHOOTING","D8",""),
    @("SWEEP","FIGHTING","D8",""),
    @("TRADEMARK_WEAPON","FIGHTING","D10",""),
    @("TRADEMARK_WEAPON","SHOOTING","D10",""),
    @("TACTICIAN","KNOWLEDGE","D6","Battle"),
    @("SOUL_DRAIN","KNOWLEDGE","D10","Arcana"),
    @("ADEPT","FAITH","D8",""),
    @("ADEPT","FIGHTING","D8",""),
    @("ASSASSIN","CLIMBING","D6",""),
    @("ASSASSIN","FIGHTING","D6",""),
    @("ASSASSIN","STEALTH","D8",""),
    @("CHAMPION","FAITH","D6",""),
    @("CHAMPION","FIGHTING","D8",""),
    @("GADGETEER","WEIRD_SCIENCE","D8",""),
    @("GADGETEER","KNOWLEDGE","D6","ANY"),
    @("GADGETEER","KNOWLEDGE","D6","ANY"),
    @("HOLY_UNHOLY_WARRIOR","FAITH","D6",""),
    @("INVESTIGATOR","INVESTIGATION","D8",""),
    @("INVESTIGATOR","STREETWISE","D8",""),
    @("MCGUYVER","REPAIR","D6",""),
    @("MCGUYVER","NOTICE","D8",""),
    @("MENTALIST","PSIONICS","D6",""),
    @("MR_FIX_IT","REPAIR","D8",""),
    @("MR_FIX_IT","WEIRD_SCIENCE","D8",""),
    @("MR_FIX_IT","KNOWLEDGE","D6","ANY"),
    @("MR_FIX_IT","KNOWLEDGE","D6","ANY"),
    @("THIEF","CLIMBING","D6",""),
    @("THIEF","LOCKPICKING","D6",""),
    @("THIEF","STEALTH","D8",""),
    @("WIZARD","KNOWLEDGE","D8","Arcana"),
    @("WIZARD","SPELLCASTING","D6",""),
    @("WOODSMAN","SURVIVAL","D8",""),
    @("WOODSMAN","TRACKING","D8",""),
    @("STRONG_WILLED","INTIMIDATION","D6",""),
    @("STRONG_WILLED","TAUNT","D6",""),
    @("MARTIAL_ARTS_MASTER","FIGHTING","D12",""),
    @("WEAPON_MASTER","FIGHTING","D12","")
)

$listRows = @(
    @("BOATING","D4"),
    @("CLIMBING","D6"),
    @("DRIVING","D8"),
    @("FAITH","D10"),
    @("FIGHTING","D12"),
    @("GAMBLING",""),
    @("HEALING",""),
    @("INTIMIDATION",""),
    @("INVESTIGATION",""),
    @("KNOWLEDGE",""),
    @("LOCKPICKING",""),
    @("NOTICE",""),
    @("PERSUASION",""),
    @("PILOTING",""),
    @("PSIONICS",""),
    @("REPAIR",""),
    @("RIDING",""),
    @("SHOOTING",""),
    @("SPELLCASTING",""),
    @("STEALTH",""),
    @("STREETWISE",""),
    @("SURVIVAL",""),
    @("SWIMMING",""),
    @("TAUNT",""),
    @("THROWING",""),
    @("TRACKING",""),
    @("WEIRD_SCIENCE","")
)

# --- wipe the old body, keep the header row ----------------------------

$ws.Range("A2:D141").ClearContents()

# Pre-intern the brand-new shared strings (new skill names + new "school"
# tags) in the same order the source workbook introduced them, via a
# scratch cell well outside the used range, so the saved sharedStrings
# table lists them in that order.
$newStrings = @("Battle","Arcana","FAITH","PSIONICS","SPELLCASTING","WEIRD_SCIENCE","ANY")
foreach ($s in $newStrings) {
    $ws.Range("Z1").Value = $s
}
$ws.Range("Z1").ClearContents()

# --- rewrite the compact edge/skill/dieType/school table ---------------

$r = 2
foreach ($row in $edgeRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne "") {
        $ws.Cells.Item($r, 4).Value = $row[3]
    }
    $r++
}

# --- rewrite the skill / dieType lookup table at rows 50-76 ------------

$r = 50
foreach ($row in $listRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne "") {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    $r++
}

# --- point the data validations at the new lookup ranges ---------------

$ws.Range("B2:B113").Validation.Delete()
$ws.Range("C2:C113").Validation.Delete()

$ws.Range("B2:B44").Validation.Add(3, 1, 1, "=`$A`$50:`$A`$76")
$ws.Range("C2:C44").Validation.Add(3, 1, 1, "=`$B`$50:`$B`$54")

# --- selection / scroll position match the saved view ------------------

$ws.Range("B44").Select()
